$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.015913605690002
$ws.Range("B1").Value = 2.028582811355591
$ws.Range("C1").Value = 3.637657403945923
$ws.Range("D1").Value = 2.070173263549805
$ws.Range("E1").Value = 0.3596626818180084
